# Update Name of Algo
# Apply updated RandomForest imputation results to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.256799999999994
$ws.Range("D6").Value = -7.919200000000003
$ws.Range("D7").Value = -7.196399999999993
$ws.Range("B8").Value = 5.470299999999998
$ws.Range("D8").Value = -8.013899999999996
$ws.Range("E11").Value = 13.23119999999999
$ws.Range("A12").Value = -22.722
$ws.Range("B12").Value = 5.416299999999999
$ws.Range("B14").Value = 8.983000000000006
$ws.Range("E14").Value = 13.23340000000001
$ws.Range("D19").Value = -8.111799999999993
$ws.Range("E19").Value = 13.47690000000001
$ws.Range("D21").Value = -7.709300000000002
$ws.Range("E21").Value = 13.8011
$ws.Range("B22").Value = 5.086500000000005
$ws.Range("D24").Value = -7.996799999999996
